$d = $word.ActiveDocument
$wdAlignRight = [int][Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphRight

# ---------------------------------------------------------------------------
# 1) Insert the student signature block (right-aligned) right after the
#    "ASSIGNMENT - WEB APPLICATION DEVELOPMENT" heading, before the first
#    horizontal-rule paragraph.
# ---------------------------------------------------------------------------
$heading = $d.Paragraphs.Item(2)
$r = $heading.Range
$r.Collapse(0)   # wdCollapseEnd
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item(3)
$p1.Range.Text = "Shruthi S V"
$p1.Style = "Body Text"
$p1.Alignment = $wdAlignRight

$r = $d.Paragraphs.Item(3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(4)
$p2.Range.Text = "2024TM93518"
$p2.Style = "Body Text"
$p2.Alignment = $wdAlignRight

$r = $d.Paragraphs.Item(4).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(5)
$p3.Range.Text = "Group 112"
$p3.Style = "Body Text"
$p3.Alignment = $wdAlignRight

# ---------------------------------------------------------------------------
# 2) Append the "Drive Video Link" / "Git hub repository link" block at the
#    very end of the document (after the Summary paragraph).
# ---------------------------------------------------------------------------

# -- blank line 1 --------------------------------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$blank1.Style = "Body Text"

# -- blank line 2 --------------------------------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$blank2.Style = "Body Text"

# -- "Drive Video Link" heading ------------------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$driveHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$driveHeading.Style = "Body Text"
$driveHeading.Range.Text = "Drive Video Link"
$driveHeading.Range.Font.Bold = 1

# -- Drive link paragraph (real hyperlink) -------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$driveLinkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$driveLinkPara.Style = "Body Text"
$driveLinkPara.Range.Font.Bold = 0
$driveLinkPara.Range.Font.Underline = 1
$driveUrl = "https://drive.google.com/file/d/13mXBkvyOYaS7sSKbWiP8eTO5iXGK8fXs/view?usp=drive_link"
$driveLinkPara.Range.InsertBefore($driveUrl)

$driveLinkParaAfter = $d.Paragraphs.Item($d.Paragraphs.Count)
$driveTextRange = $d.Range($driveLinkParaAfter.Range.Start, $driveLinkParaAfter.Range.End - 1)
$driveHyperlink = $d.Hyperlinks.Add($driveTextRange, $driveUrl, "", "", $driveUrl)
$driveHyperlink.Range.Font.Underline = 1

# -- blank line 3 --------------------------------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$blank3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$blank3.Style = "Body Text"

# -- "Git hub repository link" heading -----------------------------------
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$gitHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$gitHeading.Style = "Body Text"
$gitHeading.Range.Text = "Git hub repository link"
$gitHeading.Range.Font.Bold = 1

# -- GitHub link paragraph (plain underlined text, not a real hyperlink) --
$r = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$gitLinkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$gitLinkPara.Style = "Body Text"
$gitLinkPara.Range.Font.Bold = 0
$gitLinkPara.Range.Font.Underline = 1
$gitUrl = "https://github.com/2024tm93518-cmd/SchoolEquipmentLendingSystem"
$gitLinkPara.Range.InsertBefore($gitUrl)

Write-Output "done"
